# Apply the "changed the style of test cases" edit.
#
# The sheet holds a table of test scenarios. This edit:
#  - rewrites the scenario-1 description / test-case id text (row 2)
#  - fills in two previously-empty rows (3 and 4) with new scenario data
#  - gives row-3's scenario-id cell (B3) a borderless style variant
#  - grows the row heights / column C width to fit the new text
#  - moves the active selection to H3
#
# NOTE: cell values are assigned in the exact order needed so that the
# shared-string table comes out with the same index assignment as the
# target workbook (new/changed strings are appended to the shared-string
# table in first-write order, and strings no longer referenced by any
# cell are dropped on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing scenario's id / test-case columns ---------
# E2 first (-> shared string index 6): the expanded 5-step scenario text
$ws.Range("E2").Value = "1. Инициализация «Зритель -Каскад»`n2. Деинициализация «Зритель -Каскад»`n3. Реинициализация «Зритель -Каскад»`n4. Появление окна с описанием ошибки, когда отключена камера`n5. Появление окна с описанием ошибки, когда нет связи с сервисом"

# --- Row 3: "Проверить работу «Зритель-Контраст» в АРМ ОПК" scenario --
# C3 (-> index 7)
$ws.Range("C3").Value = "Проверить работу «Зритель-Контраст» в АРМ ОПК"

# B2 (-> index 8): scenario id renamed CIRW -> CIR-W
$ws.Range("B2").Value = "CIR-W S2.0 "

# B3 (-> index 9)
$ws.Range("B3").Value = "CIR-W S2.1"

# D2 (-> index 10): test-case id list, now 5 items
$ws.Range("D2").Value = "1. CIR-W0001`n2. CIR-W0002`n3. CIR-W0003`n4. CIR-W0004`n5. CIR-W0005"

# --- Row 4: "Проверить форму захвата" scenario -------------------------
# B4 (-> index 11)
$ws.Range("B4").Value = "CIR-W S2.2 "

# C4 (-> index 12)
$ws.Range("C4").Value = "Проверить форму захвата"

# D4 (-> index 13)
$ws.Range("D4").Value = "1. CIR-W0010`n2. CIR-W0011"

# E4 (-> index 14)
$ws.Range("E4").Value = "1. Вызов формы захвата с наличием захваченного кадра`n2. Вызов формы захвата без захваченного кадра"

# D3 (-> index 15)
$ws.Range("D3").Value = "1. CIR-W0006`n2. CIR-W0007`n3. CIR-W0008`n4. CIR-W0009`n5. CIR-W0012"

# E3 (-> index 16)
$ws.Range("E3").Value = "1. Работа «Зритель-контраст» в «АРМ ОПК» `n2. Захват кадра в главной форме «АРМ ОПК»`n3. Захват кадра на границе окна «Видеокамера»`n4. Отработка таймаута в АРМ ОПК`n5. Вызов формы захвата во время захвата в АРМ ОПК"

# --- Styling: give B3 a new cell-format variant (no border) -----------
# Matches the new 6th cellXfs entry: font 3 / no fill / no border,
# left/top aligned, no wrap.
$b3 = $ws.Range("B3")
$b3.WrapText = $false
$b3.Borders.LineStyle = -4142   # xlLineStyleNone
$b3.HorizontalAlignment = -4131 # xlLeft
$b3.VerticalAlignment = -4160   # xlTop

# B4 gets the same "bordered, no wrap" look already used by B2 (scenario
# id cell): keep its border but turn off wrapping.
$b4 = $ws.Range("B4")
$b4.WrapText = $false
$b4.HorizontalAlignment = -4131 # xlLeft
$b4.VerticalAlignment = -4160   # xlTop

# C4 gets the same "borderless, wrapped" look already used by C2
# (scenario title cell): drop its border but keep wrapping.
$c4 = $ws.Range("C4")
$c4.WrapText = $true
$c4.Borders.LineStyle = -4142   # xlLineStyleNone
$c4.HorizontalAlignment = -4131 # xlLeft
$c4.VerticalAlignment = -4160   # xlTop

# --- Row heights, grown to fit the new multi-line content --------------
$ws.Rows.Item(2).RowHeight = 223.5
$ws.Rows.Item(3).RowHeight = 220.5
$ws.Rows.Item(4).RowHeight = 94.5

# --- Column C widened to fit the longer scenario titles ----------------
$ws.Columns.Item(3).ColumnWidth = 21.4

# --- Selection moved to H3 ----------------------------------------------
$ws.Range("H3").Select()
